$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.054.90"

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.782.02"
$ws.Range("E3").Value = "  +3.16%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.23%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'243.85"
$ws.Range("E5").Value = "  +0.74%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'0.9989"
$ws.Range("E6").Value = "  +0.13%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4913"
$ws.Range("E7").Value = "  -0.28%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2678"
$ws.Range("E8").Value = "  +2.21%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06269"
$ws.Range("E9").Value = "  +0.59%  "

# Row 10 - WrappedEther
$ws.Range("D10").Value = "1.784.97"
$ws.Range("E10").Value = "  +3.37%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'16.42"
$ws.Range("E11").Value = "  +3.27%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.07025"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "'0.6277"
$ws.Range("E13").Value = "  +2.55%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'4.645"
$ws.Range("E14").Value = "  +2.92%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'80.34"
$ws.Range("E15").Value = "  +3.91%  "

# Row 16 - was Dai, now WrappedBTC (rows 16 & 17 swapped)
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "28.018.32"
$ws.Range("E16").Value = "  +4.93%  "

# Row 17 - was WrappedBTC, now Dai
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "'0.9999"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18 - BinanceUSD
$ws.Range("D18").Value = "'0.9990"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.000007219"
$ws.Range("E19").Value = "  +0.24%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  +4.61%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.010.83"
$ws.Range("E21").Value = "  +3.00%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.572"
$ws.Range("E22").Value = "  +2.69%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'8.713"
$ws.Range("E23").Value = "  +1.21%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'5.228"
$ws.Range("E24").Value = "  +2.25%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'141.94"
$ws.Range("E25").Value = "  +2.61%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'15.79"
$ws.Range("E26").Value = "  +2.66%  "

# Row 27 - LidoDAOToken
$ws.Range("D27").Value = "'1.861"
$ws.Range("E27").Value = "  +5.84%  "

# Row 28 - BitcoinCash
$ws.Range("D28").Value = "'109.54"
$ws.Range("E28").Value = "  +3.09%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +0.16%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "'4.195"
$ws.Range("E30").Value = "  +7.13%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08277"
$ws.Range("E31").Value = "  +3.58%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.798"
$ws.Range("E32").Value = "  +3.33%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.04906"
$ws.Range("E33").Value = "  +8.98%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "'1.076"
$ws.Range("E34").Value = "  +7.23%  "

# Row 35 - was HuobiToken, now ImmutableX (rows 35 & 36 swapped)
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6531"
$ws.Range("E35").Value = "  +3.99%  "

# Row 36 - was ImmutableX, now HuobiToken
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.611"
$ws.Range("E36").Value = "  +0.14%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "'0.9512"
$ws.Range("E37").Value = "  +1.95%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "'2.580"
$ws.Range("E38").Value = "  +6.78%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'2.053"
$ws.Range("E39").Value = "  +0.81%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "'5.983"
$ws.Range("E40").Value = "  +6.66%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "'0.01555"
$ws.Range("E41").Value = "  +2.49%  "

# Row 42 - PaxDollar
$ws.Range("D42").Value = "'0.9996"
$ws.Range("E42").Value = "  +0.15%  "

# Row 43 - Quant
$ws.Range("D43").Value = "'99.93"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "'0.3981"
$ws.Range("E44").Value = "  +2.90%  "

# Row 45 - Aptos
$ws.Range("D45").Value = "'7.181"
$ws.Range("E45").Value = "  +3.90%  "

# Row 46 - Algorand
$ws.Range("D46").Value = "'0.1203"
$ws.Range("E46").Value = "  +3.78%  "

# Row 47 - Cronos
$ws.Range("D47").Value = "'0.05425"
$ws.Range("E47").Value = "  +0.79%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'8.029"
$ws.Range("E48").Value = "  +2.27%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "'1.306"
$ws.Range("E49").Value = "  +5.67%  "

# Row 50 - Elrond
$ws.Range("E50").Value = "  +0.86%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'52.99"
$ws.Range("E51").Value = "  +2.30%  "
